# Generate Report for Handoff
# File "b.md.md" has had a new handoff generated: update its status and
# the associated handoff file/datetime on the Overview sheet and on each
# locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md.md" row ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 is the "b.md.md" row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-25 06:54:37"

# --- de-de sheet: row 3 is the "b.md.md" row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-25 06:54:49"
